$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 536.1739
$ws.Range("I19").Value = 504
$ws.Range("J19").Value = 550.25
$ws.Range("K19").Value = 504
$ws.Range("L19").Value = 550.25
$ws.Range("M19").Value = -329
$ws.Range("N19").Value = -900.25
$ws.Range("H40").Value = 27028422
$ws.Range("I40").Value = 1275.6875
$ws.Range("J40").Value = 47620532
$ws.Range("K40").Value = 1275.6875
$ws.Range("L40").Value = 47620532
$ws.Range("M40").Value = -1100.6875
$ws.Range("N40").Value = -47620882
$ws.Range("H69").Value = 4000
$ws.Range("J69").Value = 4000
$ws.Range("L69").Value = 12000
$ws.Range("N69").Value = -13748
$ws.Range("H72").Value = 4000
$ws.Range("J72").Value = 4000
$ws.Range("L72").Value = 36000
$ws.Range("N72").Value = -44736
$ws.Range("H113").Value = 33336398
$ws.Range("I113").Value = 55557668
$ws.Range("J113").Value = 4493.6665
$ws.Range("K113").Value = 55557668
$ws.Range("L113").Value = 4493.6665
$ws.Range("M113").Value = -55554414
$ws.Range("N113").Value = -11001.6665
$ws.Range("H116").Value = 6267.5415
$ws.Range("I116").Value = 5966.077
$ws.Range("J116").Value = 6623.8184
$ws.Range("K116").Value = 5966.077
$ws.Range("L116").Value = 6623.8184
$ws.Range("M116").Value = -2524.077
$ws.Range("N116").Value = -13507.8184
$ws.Range("H132").Value = 2214.2932
$ws.Range("I132").Value = 2048.1956
$ws.Range("J132").Value = 2851
$ws.Range("K132").Value = 6144.5868
$ws.Range("L132").Value = 8553
$ws.Range("M132").Value = -3614.5868
$ws.Range("N132").Value = -13613
$ws.Range("H135").Value = 19101.785
$ws.Range("I135").Value = 23192.31
$ws.Range("J135").Value = 2367.818
$ws.Range("K135").Value = 208730.79
$ws.Range("L135").Value = 21310.362
$ws.Range("M135").Value = -206195.79
$ws.Range("N135").Value = -26380.362

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2542368.2
$ws.Range("I32").Value = 4633.618
$ws.Range("K32").Value = 4633.618
$ws.Range("M32").Value = -4346.618
$ws.Range("H45").Value = 3572.16
$ws.Range("I45").Value = 2901.9167
$ws.Range("J45").Value = 4190.846
$ws.Range("K45").Value = 2901.9167
$ws.Range("L45").Value = 4190.846
$ws.Range("M45").Value = -2524.9167
$ws.Range("N45").Value = -4944.846
$ws.Range("H61").Value = 1341.8864
$ws.Range("I61").Value = 1323
$ws.Range("K61").Value = 1323
$ws.Range("M61").Value = -1111
$ws.Range("H74").Value = 932.53656
$ws.Range("I74").Value = 921.1
$ws.Range("K74").Value = 921.1
$ws.Range("M74").Value = -47.10000000000002
$ws.Range("H77").Value = 932.53656
$ws.Range("I77").Value = 921.1
$ws.Range("K77").Value = 4605.5
$ws.Range("M77").Value = -237.5
$ws.Range("H122").Value = 2059.0688
$ws.Range("I122").Value = 1849.32
$ws.Range("J122").Value = 3370
$ws.Range("K122").Value = 5547.96
$ws.Range("L122").Value = 10110
$ws.Range("M122").Value = -3097.96
$ws.Range("N122").Value = -15010
$ws.Range("H132").Value = 100241.04
$ws.Range("I132").Value = 132799.88
$ws.Range("J132").Value = 5069.077
$ws.Range("K132").Value = 398399.64
$ws.Range("L132").Value = 15207.231
$ws.Range("M132").Value = -395869.64
$ws.Range("N132").Value = -20267.231
$ws.Range("H136").Value = 1341.8864
$ws.Range("I136").Value = 1323
$ws.Range("K136").Value = 3969
$ws.Range("M136").Value = -1419

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 133.5
$ws.Range("I22").Value = 140.2
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 140.2
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 32.80000000000001
$ws.Range("N22").Value = -446
$ws.Range("H53").Value = 41748
$ws.Range("J53").Value = 41748
$ws.Range("L53").Value = 41748
$ws.Range("N53").Value = -42896
$ws.Range("H134").Value = 61407.3
$ws.Range("I134").Value = 74313.09
$ws.Range("K134").Value = 222939.27
$ws.Range("M134").Value = -220404.27

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1991.0526
$ws.Range("I31").Value = 1488.5714
$ws.Range("J31").Value = 3398
$ws.Range("K31").Value = 1488.5714
$ws.Range("L31").Value = 3398
$ws.Range("M31").Value = -1193.5714
$ws.Range("N31").Value = -3988
$ws.Range("H34").Value = 1991.0526
$ws.Range("I34").Value = 1488.5714
$ws.Range("J34").Value = 3398
$ws.Range("K34").Value = 1488.5714
$ws.Range("L34").Value = 3398
$ws.Range("M34").Value = -1286.5714
$ws.Range("N34").Value = -3802
$ws.Range("H81").Value = 30000
$ws.Range("I81").Value = 30000
$ws.Range("K81").Value = 30000
$ws.Range("M81").Value = -29002
$ws.Range("H84").Value = 30000
$ws.Range("I84").Value = 30000
$ws.Range("K84").Value = 90000
$ws.Range("M84").Value = -85008
$ws.Range("H86").Value = 62510080
$ws.Range("I86").Value = 90921130
$ws.Range("J86").Value = 5761.6
$ws.Range("K86").Value = 90921130
$ws.Range("L86").Value = 5761.6
$ws.Range("M86").Value = -90920007
$ws.Range("N86").Value = -8007.6
$ws.Range("H89").Value = 62510080
$ws.Range("I89").Value = 90921130
$ws.Range("J89").Value = 5761.6
$ws.Range("K89").Value = 454605650
$ws.Range("L89").Value = 28808
$ws.Range("M89").Value = -454600034
$ws.Range("N89").Value = -40040
$ws.Range("H132").Value = 2254.375
$ws.Range("I132").Value = 1820.7778
$ws.Range("J132").Value = 3555.1667
$ws.Range("K132").Value = 5462.3334
$ws.Range("L132").Value = 10665.5001
$ws.Range("M132").Value = -2932.3334
$ws.Range("N132").Value = -15725.5001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1293.2646
$ws.Range("I129").Value = 690.8182
$ws.Range("J129").Value = 1581.3914
$ws.Range("K129").Value = 2072.4546
$ws.Range("L129").Value = 4744.174199999999
$ws.Range("M129").Value = 2927.5454
$ws.Range("N129").Value = -14744.1742

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 657.4516
$ws.Range("I107").Value = 467.73685
$ws.Range("J107").Value = 957.8333
$ws.Range("K107").Value = 467.73685
$ws.Range("L107").Value = 957.8333
$ws.Range("M107").Value = 1452.26315
$ws.Range("N107").Value = -4797.8333
$ws.Range("H113").Value = 1275.2693
$ws.Range("I113").Value = 1097.625
$ws.Range("J113").Value = 1559.5
$ws.Range("K113").Value = 1097.625
$ws.Range("L113").Value = 1559.5
$ws.Range("M113").Value = 1072.375
$ws.Range("N113").Value = -5899.5
$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 3500
$ws.Range("J122").Value = 2166.6667
$ws.Range("K122").Value = 10500
$ws.Range("L122").Value = 6500.000100000001
$ws.Range("M122").Value = -8050
$ws.Range("N122").Value = -11400.0001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11113344
$ws.Range("I7").Value = 20001660
$ws.Range("J7").Value = 2949.75
$ws.Range("K7").Value = 20001660
$ws.Range("L7").Value = 2949.75
$ws.Range("M7").Value = -20001548
$ws.Range("N7").Value = -3173.75
$ws.Range("H122").Value = 2288.0908
$ws.Range("I122").Value = 1631.0834
$ws.Range("K122").Value = 4893.2502
$ws.Range("M122").Value = -2443.2502
$ws.Range("H126").Value = 11113344
$ws.Range("I126").Value = 20001660
$ws.Range("J126").Value = 2949.75
$ws.Range("K126").Value = 60004980
$ws.Range("L126").Value = 8849.25
$ws.Range("M126").Value = -60002510
$ws.Range("N126").Value = -13789.25

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 30271
$ws.Range("J69").Value = 30271
$ws.Range("L69").Value = 30271
$ws.Range("N69").Value = -31769
$ws.Range("H72").Value = 30271
$ws.Range("J72").Value = 30271
$ws.Range("L72").Value = 90813
$ws.Range("N72").Value = -98301
$ws.Range("H75").Value = 27644.285
$ws.Range("J75").Value = 27644.285
$ws.Range("L75").Value = 27644.285
$ws.Range("N75").Value = -29516.285
$ws.Range("H78").Value = 27644.285
$ws.Range("J78").Value = 27644.285
$ws.Range("L78").Value = 82932.855
$ws.Range("N78").Value = -92292.855
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 6409.2856
$ws.Range("I81").Value = 10576.363
$ws.Range("J81").Value = 1825.5
$ws.Range("K81").Value = 21152.726
$ws.Range("L81").Value = 3651
$ws.Range("M81").Value = -20091.726
$ws.Range("N81").Value = -5773
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 6409.2856
$ws.Range("I84").Value = 10576.363
$ws.Range("J84").Value = 1825.5
$ws.Range("K84").Value = 105763.63
$ws.Range("L84").Value = 18255
$ws.Range("M84").Value = -100459.63
$ws.Range("N84").Value = -28863
$ws.Range("H87").Value = 32000
$ws.Range("J87").Value = 32000
$ws.Range("L87").Value = 32000
$ws.Range("M87").Value = -34496
$ws.Range("H90").Value = 32000
$ws.Range("J90").Value = 32000
$ws.Range("L90").Value = 96000
$ws.Range("N90").Value = -108480
$ws.Range("H132").Value = 2099.7856
$ws.Range("I132").Value = 2215.2327
$ws.Range("J132").Value = 1717.9231
$ws.Range("K132").Value = 6645.6981
$ws.Range("L132").Value = 5153.7693
$ws.Range("M132").Value = -4115.6981
$ws.Range("N132").Value = -10213.7693
